$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the number format of A9 (was date-only, should be the standard datetime format
# used by the rest of the column).
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append new row 10 with the latest transfer record.
$ws.Range("A10").Value = 45856
$ws.Range("A10").NumberFormat = "YYYY-MM-DD"

$ws.Range("B10").Value = "QWE1234"
$ws.Range("C10").Value = "SAMUEL"
$ws.Range("D10").Value = "2025-07-18 13:54:39"
$ws.Range("E10").Value = "2025-07-18 13:54:40"
$ws.Range("F10").Value = "2025-07-18 13:54:40"
$ws.Range("G10").Value = "2025-07-18 13:54:42"
$ws.Range("H10").Value = "2025-07-18 13:54:42"
$ws.Range("I10").Value = "2025-07-18 13:54:44"
$ws.Range("J10").Value = "2025-07-18 13:54:45"
$ws.Range("K10").Value = "0:00:02"
$ws.Range("L10").Value = "0:00:01"
$ws.Range("M10").Value = "0:00:06"
$ws.Range("N10").Value = ""
$ws.Range("O10").Value = "2025-07-18 13:54:47"
$ws.Range("P10").Value = "2025-07-18 13:54:48"
$ws.Range("Q10").Value = "2025-07-18 13:54:49"
$ws.Range("R10").Value = "2025-07-18 13:54:50"
$ws.Range("S10").Value = "0:00:01"
$ws.Range("T10").Value = "0:00:01"
$ws.Range("U10").Value = "0:00:04"
$ws.Range("V10").Value = "0:00:01"
$ws.Range("W10").Value = "2025-07-18 13:54:46"
